# Apply the edits described by the diff: update a set of single-column
# table cells with new statistic values (some cells collapse several
# tab-separated runs into one new value).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row (1-based) -> new cell text
$changes = @{
    1  = "0M"
    2  = "0M"
    3  = "0M"
    4  = "561"
    5  = "0.00001"
    6  = "0.00081"
    7  = "0.00018"
    8  = "0.00005"
    9  = "0.00030"
    10 = "0.00041"
    11 = "0.00050"
    12 = "0.11981"
    44 = "99.87"
    45 = "0.12"
    46 = "90"
}

foreach ($rowIndex in $changes.Keys) {
    $cell = $t.Cell($rowIndex, 1)
    $cell.Range.Text = $changes[$rowIndex]
}
